$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "корпорация" (corporation) column header and sample entity data
$ws.Range("N1").Value = "корпорация"
$ws.Range("N2").Value = "Microsoft"
$ws.Range("N3").Value = "Tesla"
$ws.Range("N4").Value = "SpaceX"

# Update the active selection to reflect where editing left off
$ws.Range("N5").Select()
